# ---------------------------------------------------------------------------
# Commit "commit 2 from abhishek":
#   - consolidate the per-suite TestNG parameters (previously their own sheet,
#     "suiteLevelParametersTestNG") into new columns on the main details sheet
#   - rename sheets: projectStructureDetails -> projectLevelDetails,
#                    testCases               -> testCasesTestNG
#   - delete the now-redundant suiteLevelParametersTestNG sheet
#   - rename the "comma separated" columns/headers on testCasesTestNG to the
#     shorter "_csv" naming, and store their example values as actual
#     multi-line (wrapped) text instead of a single comma-separated line
#   - update the tcquery sample query text to reference the renamed sheet
#   - move the active sheet/selection around a bit (cosmetic)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper: paste the number/alignment/fill formatting of $src onto $dst ---
function Copy-Format($src, $dst) {
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 1) projectStructureDetails: add the two new "rolled up" TestNG columns
# ---------------------------------------------------------------------------
$wsDetails = $wb.Worksheets.Item("projectStructureDetails")
$wsCasesForFmt = $wb.Worksheets.Item("testCases")

# Headers (row 1) - K gets the same look as the rest of the header row,
# L keeps the plain yellow fill (no wrap), which on this workbook is the
# style used by testCases' own header row (style "1": fill, no wrap).
$wsDetails.Cells.Item(1, 11).Value = "testNg_listenerClasses_csv"
Copy-Format $wsDetails.Cells.Item(1, 1) $wsDetails.Cells.Item(1, 11)

$wsDetails.Cells.Item(1, 12).Value = "testNG_suite_level_paramName:paramValue_csv"
Copy-Format $wsCasesForFmt.Cells.Item(1, 1) $wsDetails.Cells.Item(1, 12)

# Sample data (row 2) - wrapped, multi-line text pulled in from the old
# suiteLevelParametersTestNG sheet.
$wsDetails.Cells.Item(2, 11).Value = "listenerclass1,`nlistenerclass2"
$wsDetails.Cells.Item(2, 11).WrapText = $true

$wsDetails.Cells.Item(2, 12).Value = "suiteparamname1:suiteparamvalue1,`nsuiteparamname2:suiteparamvalue2`n"
$wsDetails.Cells.Item(2, 12).WrapText = $true

$wsDetails.Columns.Item(11).ColumnWidth = 28.7109375
$wsDetails.Columns.Item(12).ColumnWidth = 54.140625

$wsDetails.Rows.Item(2).RowHeight = 45

# ---------------------------------------------------------------------------
# 2) testCases: rename the csv columns + store wrapped multi-line samples
# ---------------------------------------------------------------------------
$wsCases = $wb.Worksheets.Item("testCases")

$wsCases.Cells.Item(1, 4).Value = "testClasses_csv"
$wsCases.Cells.Item(1, 6).Value = "test_paramName:paramValue_csv"

$wsCases.Cells.Item(2, 4).Value = "testclass1,`ntestclass2"
$wsCases.Cells.Item(2, 4).WrapText = $true
$wsCases.Cells.Item(2, 6).Value = "paramname1:paramvalue1,`nparamname2:paramvalue2"
$wsCases.Cells.Item(2, 6).WrapText = $true

$wsCases.Cells.Item(3, 4).Value = "testclass1,`ntestclass3"
$wsCases.Cells.Item(3, 4).WrapText = $true
$wsCases.Cells.Item(3, 6).Value = "paramname11:paramvalue11,`nparamname12:paramvalue12"
$wsCases.Cells.Item(3, 6).WrapText = $true

$wsCases.Cells.Item(4, 4).Value = "testclass1,`ntestclass4"
$wsCases.Cells.Item(4, 4).WrapText = $true
$wsCases.Cells.Item(4, 6).Value = "paramname11:paramvalue11,`nparamname12:paramvalue13"
$wsCases.Cells.Item(4, 6).WrapText = $true

$wsCases.Cells.Item(5, 4).Value = "testclass1,`ntestclass5"
$wsCases.Cells.Item(5, 4).WrapText = $true
$wsCases.Cells.Item(5, 6).Value = "paramname11:paramvalue11,`nparamname12:paramvalue14"
$wsCases.Cells.Item(5, 6).WrapText = $true

$wsCases.Cells.Item(6, 4).Value = "testclass1"
$wsCases.Cells.Item(6, 6).Value = "paramname11:paramvalue11,`nparamname12:paramvalue15"
$wsCases.Cells.Item(6, 6).WrapText = $true

$wsCases.Rows.Item(2).RowHeight = 30
$wsCases.Rows.Item(3).RowHeight = 30
$wsCases.Rows.Item(4).RowHeight = 30
$wsCases.Rows.Item(5).RowHeight = 30
$wsCases.Rows.Item(6).RowHeight = 30

# ---------------------------------------------------------------------------
# 3) queries: tcquery's sample SQL now points at the renamed sheet
# ---------------------------------------------------------------------------
$wsQueries = $wb.Worksheets.Item("queries")
$wsQueries.Cells.Item(3, 2).Value = "select * from testCasesTestNG where module='module1'"

# ---------------------------------------------------------------------------
# 4) drop the now-redundant suite-level-parameters sheet
# ---------------------------------------------------------------------------
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("suiteLevelParametersTestNG").Delete()
$excel.DisplayAlerts = $true

# ---------------------------------------------------------------------------
# 5) rename sheets
# ---------------------------------------------------------------------------
$wsDetails.Name = "projectLevelDetails"
$wsCases.Name = "testCasesTestNG"

# ---------------------------------------------------------------------------
# 6) restore per-sheet selections (cosmetic, matches the saved view state)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("dependenciesMaven").Activate()
$wb.Worksheets.Item("dependenciesMaven").Range("C9").Select()

$wb.Worksheets.Item("projectSkeletons").Activate()
$wb.Worksheets.Item("projectSkeletons").Range("B5").Select()

$wb.Worksheets.Item("queries").Activate()
$wb.Worksheets.Item("queries").Range("B15").Select()

$wsDetails.Activate()
$wsDetails.Range("I15").Select()

# testCasesTestNG is the sheet left active/selected in the saved workbook.
$wsCases.Activate()
$wsCases.Range("C16").Select()
